$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 255.95454
$ws.Range("I28").Value = 256.6
$ws.Range("J28").Value = 249.5
$ws.Range("K28").Value = 256.6
$ws.Range("L28").Value = 249.5
$ws.Range("M28").Value = 228.4
$ws.Range("N28").Value = -1219.5
# Row 33
$ws.Range("H33").Value = 123
$ws.Range("I33").Value = 128.4
$ws.Range("K33").Value = 128.4
$ws.Range("M33").Value = 100.6
# Row 43
$ws.Range("H43").Value = 12554.454
$ws.Range("I43").Value = 13500
$ws.Range("J43").Value = 12459.9
$ws.Range("K43").Value = 13500
$ws.Range("L43").Value = 12459.9
$ws.Range("M43").Value = -13431
$ws.Range("N43").Value = -12597.9
# Row 62
$ws.Range("H62").Value = 8995
$ws.Range("I62").Value = 5995
$ws.Range("J62").Value = 11995
$ws.Range("K62").Value = 5995
$ws.Range("L62").Value = 11995
$ws.Range("M62").Value = -5371
$ws.Range("N62").Value = -13243
# Row 65
$ws.Range("H65").Value = 8995
$ws.Range("I65").Value = 5995
$ws.Range("J65").Value = 11995
$ws.Range("K65").Value = 29975
$ws.Range("L65").Value = 59975
$ws.Range("M65").Value = -26855
$ws.Range("N65").Value = -66215
# Row 74
$ws.Range("H74").Value = 3357.6843
$ws.Range("J74").Value = 5333.3335
$ws.Range("L74").Value = 5333.3335
$ws.Range("N74").Value = -7205.3335
# Row 77
$ws.Range("H77").Value = 3357.6843
$ws.Range("J77").Value = 5333.3335
$ws.Range("L77").Value = 26666.6675
$ws.Range("N77").Value = -36026.6675
# Row 107
$ws.Range("H107").Value = 34153.87
$ws.Range("I107").Value = 640.7692
$ws.Range("K107").Value = 640.7692
$ws.Range("M107").Value = 1279.2308
# Row 111
$ws.Range("H111").Value = 266.66666
$ws.Range("I111").Value = 266.66666
$ws.Range("K111").Value = 799.9999799999999
$ws.Range("M111").Value = 2267.00002
# Row 129
$ws.Range("H129").Value = 1882.9412
$ws.Range("I129").Value = 810.1818
$ws.Range("K129").Value = 2430.5454
$ws.Range("M129").Value = 2569.4546
# Row 137
$ws.Range("H137").Value = 993.24
$ws.Range("I137").Value = 942.64703
$ws.Range("K137").Value = 2827.94109
$ws.Range("M137").Value = -277.9410899999998
# Row 141
$ws.Range("H141").Value = 2221.0715
$ws.Range("I141").Value = 2221.0715
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6663.2145
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1483.2145
$ws.Range("N141").ClearContents()

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
# Row 32
$ws.Range("H32").Value = 2916.6326
$ws.Range("I32").Value = 2261.9092
$ws.Range("K32").Value = 2261.9092
$ws.Range("M32").Value = -1974.9092
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
# Row 37
$ws.Range("H37").Value = 5021950.5
$ws.Range("I37").Value = 16681545
$ws.Range("J37").Value = 24981.428
$ws.Range("K37").Value = 16681545
$ws.Range("L37").Value = 24981.428
$ws.Range("M37").Value = -16681272
$ws.Range("N37").Value = -25527.428
# Row 45
$ws.Range("H45").Value = 8871.4
$ws.Range("I45").Value = 11864.637
$ws.Range("J45").Value = 5213
$ws.Range("K45").Value = 11864.637
$ws.Range("L45").Value = 5213
$ws.Range("M45").Value = -11487.637
$ws.Range("N45").Value = -5967
# Row 61
$ws.Range("H61").Value = 2160.6667
$ws.Range("I61").Value = 2138
$ws.Range("J61").Value = 2172
$ws.Range("K61").Value = 2138
$ws.Range("L61").Value = 2172
$ws.Range("M61").Value = -1926
$ws.Range("N61").Value = -2596
# Row 102
$ws.Range("H102").Value = 1957.75
$ws.Range("I102").Value = 1951.2333
$ws.Range("J102").Value = 2055.5
$ws.Range("K102").Value = 1951.2333
$ws.Range("L102").Value = 2055.5
$ws.Range("M102").Value = -329.2333000000001
$ws.Range("N102").Value = -5299.5
# Row 122
$ws.Range("H122").Value = 2580.1765
$ws.Range("I122").Value = 2546.3845
$ws.Range("K122").Value = 7639.1535
$ws.Range("M122").Value = -5189.1535
# Row 136
$ws.Range("H136").Value = 2160.6667
$ws.Range("I136").Value = 2138
$ws.Range("J136").Value = 2172
$ws.Range("K136").Value = 6414
$ws.Range("L136").Value = 6516
$ws.Range("M136").Value = -3864
$ws.Range("N136").Value = -11616

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4062
$ws.Range("I16").Value = 5170.3335
$ws.Range("K16").Value = 5170.3335
$ws.Range("M16").Value = -4883.3335
# Row 31
$ws.Range("H31").Value = 10148.583
$ws.Range("I31").Value = 2988.3428
$ws.Range("J31").Value = 29426.154
$ws.Range("K31").Value = 2988.3428
$ws.Range("L31").Value = 29426.154
$ws.Range("M31").Value = -2693.3428
$ws.Range("N31").Value = -30016.154
# Row 34
$ws.Range("H34").Value = 10148.583
$ws.Range("I34").Value = 2988.3428
$ws.Range("J34").Value = 29426.154
$ws.Range("K34").Value = 2988.3428
$ws.Range("L34").Value = 29426.154
$ws.Range("M34").Value = -2786.3428
$ws.Range("N34").Value = -29830.154
# Row 62
$ws.Range("H62").Value = 3499.5
$ws.Range("I62").Value = 3332.6667
$ws.Range("K62").Value = 3332.6667
$ws.Range("M62").Value = -2708.6667
# Row 65
$ws.Range("H65").Value = 3499.5
$ws.Range("I65").Value = 3332.6667
$ws.Range("K65").Value = 16663.3335
$ws.Range("M65").Value = -13543.3335
# Row 113
$ws.Range("H113").Value = 4062
$ws.Range("I113").Value = 5170.3335
$ws.Range("K113").Value = 5170.3335
$ws.Range("M113").Value = -3000.3335

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 1667278.5
$ws.Range("I11").Value = 2500139
$ws.Range("J11").Value = 1557.6666
$ws.Range("K11").Value = 7500417
$ws.Range("L11").Value = 4672.9998
$ws.Range("M11").Value = -7500277
$ws.Range("N11").Value = -4952.9998
# Row 48
$ws.Range("H48").Value = 1333.3334
$ws.Range("J48").Value = 1400
$ws.Range("L48").Value = 4200
$ws.Range("N48").Value = -4700
# Row 52
$ws.Range("H52").Value = 1409.5
$ws.Range("J52").Value = 1409.5
$ws.Range("L52").Value = 4228.5
$ws.Range("N52").Value = -4760.5
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
# Row 81
$ws.Range("H81").Value = 35723536
$ws.Range("I81").Value = 7514.3335
$ws.Range("J81").Value = 62510548
$ws.Range("K81").Value = 22543.0005
$ws.Range("L81").Value = 187531644
$ws.Range("M81").Value = -21420.0005
$ws.Range("N81").Value = -187533890
# Row 84
$ws.Range("H84").Value = 35723536
$ws.Range("I84").Value = 7514.3335
$ws.Range("J84").Value = 62510548
$ws.Range("K84").Value = 67629.0015
$ws.Range("L84").Value = 562594932
$ws.Range("M84").Value = -62013.0015
$ws.Range("N84").Value = -562606164
# Row 122
$ws.Range("H122").Value = 1211.2
$ws.Range("J122").Value = 1211.2
$ws.Range("L122").Value = 10900.8
$ws.Range("N122").Value = -15800.8
# Row 124
$ws.Range("H124").Value = 1404.5
$ws.Range("I124").Value = 1404.5
$ws.Range("K124").Value = 4213.5
$ws.Range("M124").Value = 696.5
# Row 129
$ws.Range("H129").Value = 120062.35
$ws.Range("J129").Value = 3219.3635
$ws.Range("L129").Value = 9658.0905
$ws.Range("N129").Value = -19658.0905
# Row 131
$ws.Range("H131").Value = 1868.9
$ws.Range("J131").Value = 1665.5714
$ws.Range("L131").Value = 4996.7142
$ws.Range("N131").Value = -15076.7142

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 8266.875
$ws.Range("I107").Value = 445.33334
$ws.Range("K107").Value = 445.33334
$ws.Range("M107").Value = 1474.66666
# Row 113
$ws.Range("H113").Value = 3435.8333
$ws.Range("I113").Value = 3112.5
$ws.Range("J113").Value = 3888.5
$ws.Range("K113").Value = 3112.5
$ws.Range("L113").Value = 3888.5
$ws.Range("M113").Value = -942.5
$ws.Range("N113").Value = -8228.5

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 72540.57000000001
$ws.Range("I16").Value = 1142.5454
$ws.Range("J16").Value = 334333.34
$ws.Range("K16").Value = 1142.5454
$ws.Range("L16").Value = 334333.34
$ws.Range("M16").Value = -972.5454
$ws.Range("N16").Value = -334673.34
# Row 22
$ws.Range("H22").Value = 756
$ws.Range("I22").Value = 695
$ws.Range("K22").Value = 695
$ws.Range("M22").Value = -400
# Row 27
$ws.Range("H27").Value = 756
$ws.Range("I27").Value = 695
$ws.Range("K27").Value = 695
$ws.Range("M27").Value = -588
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
# Row 93
$ws.Range("H93").Value = 10562.757
$ws.Range("I93").Value = 1630.1818
$ws.Range("J93").Value = 84256.5
$ws.Range("K93").Value = 1630.1818
$ws.Range("L93").Value = 84256.5
$ws.Range("M93").Value = -382.1818000000001
$ws.Range("N93").Value = -86752.5
# Row 136
$ws.Range("H136").Value = 4724.1113
$ws.Range("I136").Value = 3498
$ws.Range("J136").Value = 5705
$ws.Range("K136").Value = 10494
$ws.Range("L136").Value = 17115
$ws.Range("M136").Value = -7944
$ws.Range("N136").Value = -22215

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 2120.9473
$ws.Range("I100").Value = 1212.5
$ws.Range("K100").Value = 2425
$ws.Range("M100").Value = -1884
# Row 124
$ws.Range("H124").Value = 44073.5
$ws.Range("J124").Value = 44073.5
$ws.Range("L124").Value = 44073.5
$ws.Range("N124").Value = -53893.5
# Row 127
$ws.Range("H127").Value = 129992.664
$ws.Range("J127").Value = 129992.664
$ws.Range("L127").Value = 129992.664
$ws.Range("N127").Value = -139912.664
